$wb = $excel.ActiveWorkbook

# --- About sheet: add Wyoming label + date in row 1 ---
$aboutWs = $wb.Worksheets.Item("About")
$aboutWs.Range("B1").Value = "Wyoming"
$aboutWs.Range("C1").Value = [DateTime]"2024-03-12"
$aboutWs.Range("C1").NumberFormat = "mm-dd-yy"

# --- RACP sheet: wrap the formula in ROUND(...,0) ---
$racpWs = $wb.Worksheets.Item("RACP")
$racpWs.Range("B2").Formula = "=ROUND(100/About!A11,0)"
